$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "MEC-1B-T. M. Metalicos"
$ws.Range("C3").Value = "MEC-1B-T. M. Metalicos"
$ws.Range("C4").Value = "-"
$ws.Range("C6").Value = "-"

$ws.Range("E10").Value = "MEC-2A-Metalografia"

$ws.Range("D11").Value = "MEC-1A-T. M. Metalicos"
$ws.Range("E11").Value = "MEC-2A-Metalografia"

$ws.Range("D12").Value = "MEC-1A-T. M. Metalicos"
$ws.Range("E12").Value = "MEC-2A-Metalografia"
$ws.Range("F12").Value = "-"

$ws.Range("E14").Value = "MEC-2A-Metalografia"
$ws.Range("F14").Value = "-"

$ws.Range("B15").Value = "-"
$ws.Range("F15").Value = "-"

$ws.Range("B16").Value = "-"
$ws.Range("F16").Value = "-"
